$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13
$ws1.Range("F9").Value = 63
$ws1.Range("F10").Value = 459

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F9").Value = 63
$ws4.Range("F10").Value = 459
